$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows before row 38 (these become the new "Risk Management Plan"
#    and "Unit Test Documentation" rows). Everything previously at rows 38+ shifts down by 2.
$ws.Rows("38:39").Insert()

# 2. Grow the table so it covers the two new rows (A1:N47 -> A1:N49)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N49"))

# 3. Populate the two newly inserted rows.
$ws.Range("A38").Value = "Risk Management Plan"
$ws.Range("B38").Value = "Zoheb"

$ws.Range("A39").Value = "Unit Test Documentation"
$ws.Range("B39").Value = "Talhah"

# Match formatting used by the other sub rows directly above (rows 36/37 area)
$ws.Range("B38:B39").Style = "Bad"
$ws.Range("K38:N39").Style = "Bad"

# 4. Update the two "5/8 (Thursday)" dates that changed to "5/6 (Tuesday)" / "5/6 - 5/8"
$ws.Range("N31").Value = "5/6 (Tuesday)"
$ws.Range("N44").Value = "5/6 (Tuesday)"
$ws.Range("N45").Value = "5/6 - 5/8"

# 5. Leave the cursor where the author left it.
$ws.Range("I42").Select()
